$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.389.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6326"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07590"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2972"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "2.528.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +36.50%  "
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.587.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +22.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.982"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6859"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009928"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.170"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.430.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "231.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.607"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.453"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.473"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05816"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.262"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.126"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.023"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.869"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.159"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.568.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +27.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.249.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.790"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9049"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.081"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.290"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.183"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4013"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.692"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1124"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "
